# Add the METABOLIGHTS_METABOLOMICS "ER list" sheet (empty template) as the
# last sheet in the workbook, and make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

$data = @(
    ,(@('','TermSourceRef','Ontology','TAN','Content type (validation)','Notes during templating','Target term','Instruction','Requirement (m/o/n)','Value (cv/s/d)','Additional information','Review comments'))
    ,(@('Source Name','','','','','','','','','','',''))
    ,(@('Sample Name','','','','','','','','','','',''))
    ,(@('Parameter [MS sample post-extraction]','NFDI4PSO:0000043','NFDI4PSO','http://purl.obolibrary.org/obo/NFDI4PSO_0000043','','','','','','','',''))
    ,(@('Parameter [MS sample resuspension]','NFDI4PSO:0000044','NFDI4PSO','http://purl.obolibrary.org/obo/NFDI4PSO_0000044','','','','','','','',''))
    ,(@('Parameter [MS sample type]','NFDI4PSO:0000045','NFDI4PSO','http://purl.obolibrary.org/obo/NFDI4PSO_0000045','','','','','','','',''))
    ,(@('Parameter [MS derivatization]','NFDI4PSO:0000052','NFDI4PSO','http://purl.obolibrary.org/obo/NFDI4PSO_0000052','','','','','','','',''))
    ,(@('Parameter [label]','NFDI4PSO:0000079','NFDI4PSO','http://purl.obolibrary.org/obo/NFDI4PSO_0000079','','','','','','','',''))
    ,(@('Parameter [Chromatography instrument model]','NFDI4PSO:0000046','NFDI4PSO','http://purl.obolibrary.org/obo/NFDI4PSO_0000046','','','','','','','',''))
    ,(@('Parameter [Chromatography autosampler model]','NFDI4PSO:0000047','NFDI4PSO','http://purl.obolibrary.org/obo/NFDI4PSO_0000047','','','','','','','',''))
    ,(@('Parameter [Chromatography column type]','NFDI4PSO:0000053','NFDI4PSO','http://purl.obolibrary.org/obo/NFDI4PSO_0000053','','','','','','','',''))
    ,(@('Parameter [Chromatography column model]','NFDI4PSO:0000048','NFDI4PSO','http://purl.obolibrary.org/obo/NFDI4PSO_0000048','','','','','','','',''))
    ,(@('Parameter [Chromatography guard column model]','NFDI4PSO:0000049','NFDI4PSO','http://purl.obolibrary.org/obo/NFDI4PSO_0000049','','','','','','','',''))
    ,(@('Parameter [mobile phase]','NFDI4PSO:0000080','NFDI4PSO','http://purl.obolibrary.org/obo/NFDI4PSO_0000080','','','','','','','',''))
    ,(@('Parameter [elution]','NFDI4PSO:0000081','NFDI4PSO','http://purl.obolibrary.org/obo/NFDI4PSO_0000081','','','','','','','',''))
    ,(@('Parameter [scan polarity]','MS:1000465','MS','http://purl.obolibrary.org/obo/MS_1000465','','','','','','','',''))
    ,(@('Parameter [scan window lower limit]','MS:1000501','MS','http://purl.obolibrary.org/obo/MS_1000501','','','','','','','',''))
    ,(@('Parameter [scan window upper limit]','MS:1000500','MS','http://purl.obolibrary.org/obo/MS_1000500','','','','','','','',''))
    ,(@('Parameter [instrument model]','MS:1000031','MS','http://purl.obolibrary.org/obo/MS_1000031','','','','','','','',''))
    ,(@('Parameter [ionization type]','MS:1000008','MS','http://purl.obolibrary.org/obo/MS_1000008','','','','','','','',''))
    ,(@('Parameter [mass analyzer type]','MS:1000443','MS','http://purl.obolibrary.org/obo/MS_1000443','','','','','','','',''))
    ,(@('Parameter [detector type]','MS:1000026','MS','http://purl.obolibrary.org/obo/MS_1000026','','','','','','','',''))
    ,(@('Parameter [detector acquisition mode]','MS:1000027','MS','http://purl.obolibrary.org/obo/MS_1000027','','','','','','','',''))
)

$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "METABOLIGHTS_METABOLOMICS"

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $newSheet.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# The new sheet becomes the selected / active tab.
$newSheet.Select()
